$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $range = $ws.Range($cellRef)
    $escaped = $value.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

Set-TextValue 'D2' '25.759.32'
Set-TextValue 'E2' '  -0.33%  '
Set-TextValue 'D3' '1.633.15'
Set-TextValue 'E3' '  -0.18%  '
Set-TextValue 'E4' '  -0.18%  '
Set-TextValue 'D5' '214.90'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'E6' '  -0.84%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'E9' '  -1.12%  '
Set-TextValue 'D10' '19.67'
Set-TextValue 'E10' '  -3.45%  '
Set-TextValue 'D11' '0.0786'
Set-TextValue 'E11' '  +0.83%  '
Set-TextValue 'D12' '4.25'
Set-TextValue 'E12' '  -0.13%  '
Set-TextValue 'D13' '1.857.75'
Set-TextValue 'E13' '  -0.22%  '
Set-TextValue 'D14' '1.633.15'
Set-TextValue 'E14' '  -0.37%  '
Set-TextValue 'D15' '0.557'
Set-TextValue 'E15' '  -0.38%  '
Set-TextValue 'E16' '  -0.46%  '
Set-TextValue 'D17' '62.71'
Set-TextValue 'E17' '  -0.89%  '
Set-TextValue 'D18' '25.753.33'
Set-TextValue 'E19' '  -0.20%  '
Set-TextValue 'E20' '  +1.53%  '
Set-TextValue 'D21' '194.06'
Set-TextValue 'E21' '  +0.76%  '
Set-TextValue 'D22' '9.94'
Set-TextValue 'E22' '  -0.12%  '
Set-TextValue 'E23' '  +1.82%  '
Set-TextValue 'E24' '  -0.14%  '
Set-TextValue 'E25' '  +4.27%  '
Set-TextValue 'D26' '142.86'
Set-TextValue 'E26' '  +2.59%  '
Set-TextValue 'D27' '0.122'
Set-TextValue 'E27' '  -0.73%  '
Set-TextValue 'E28' '  +0.54%  '
Set-TextValue 'E29' '  +0.06%  '
Set-TextValue 'E30' '  -0.11%  '
Set-TextValue 'E31' '  -0.70%  '
Set-TextValue 'D32' '3.33'
Set-TextValue 'E32' '  +0.51%  '
Set-TextValue 'E33' '  -0.51%  '
Set-TextValue 'E34' '  +0.47%  '
Set-TextValue 'E35' '  +0.01%  '
Set-TextValue 'E36' '  -0.15%  '
Set-TextValue 'D37' '1.128.54'
Set-TextValue 'E37' '  -0.48%  '
Set-TextValue 'E38' '  -1.54%  '
Set-TextValue 'D39' '0.547'
Set-TextValue 'E39' '  -1.97%  '
Set-TextValue 'E42' '  +2.03%  '
Set-TextValue 'D43' '99.88'
Set-TextValue 'E43' '  +0.75%  '
Set-TextValue 'E44' '  +0.57%  '
Set-TextValue 'D45' '1.767.36'
Set-TextValue 'E45' '  -0.35%  '
Set-TextValue 'D46' '0.0₆0112'
Set-TextValue 'E46' '  -1.64%  '
Set-TextValue 'D47' '55.05'
Set-TextValue 'E47' '  -1.13%  '
Set-TextValue 'D48' '0.417'
Set-TextValue 'E48' '  -2.32%  '
Set-TextValue 'E49' '  +0.14%  '
Set-TextValue 'B50' 'SynthetixNetwork'
Set-TextValue 'C50' 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue 'D50' '2.34'
Set-TextValue 'E50' '  +2.81%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.56'
Set-TextValue 'E51' '  -3.23%  '

$excel.CutCopyMode = 0

Write-Output "Done applying changes"